# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (rows 3-7): refresh this week's numbers ---

# Row 3: Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2
$ws.Range("C3").Value = 1246
$ws.Range("D3").Value = 53.9

# Row 4: Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2
$ws.Range("C4").Value = 3271

# Row 5: Intel(R) Wi-Fi 6E AX211 160MHz - 22.200.2.1
$ws.Range("D5").Value = 96.6

# Row 6: driver version rolled from AX210 23.90.0.2 -> AX211 23.30.0.6
$ws.Range("A6").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.30.0.6"
$ws.Range("C6").Value = 6

# Row 7: Totals
$ws.Range("C7").Value = 4862

# --- "Good Drivers" table (rows 15-31): drop retired AX210 driver rows ---
# these adapter/driver rows no longer qualify this week and are removed,
# shifting the remaining rows up. Deleting from the bottom up keeps the
# row numbers of the not-yet-deleted rows stable.
$rowsToRemove = @(15, 16, 18, 19, 23, 24, 25, 27)
$sortedDesc = $rowsToRemove | Sort-Object -Descending
foreach ($r in $sortedDesc) {
    $ws.Rows($r).Delete()
}
